$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Every existing row's "Förändrad" (Changed) date in column C moves
#    from 45182 (2023-09-13) to 45184 (2023-09-15), for data rows 2..505.
$ws.Range("C2:C505").Value = 45184

# 2) The previously-last data row (505) picks up an explicit row height,
#    matching every other data row in the sheet.
$ws.Rows.Item(505).RowHeight = 15

# 3) Four brand-new cleavings/notifications are appended as rows 506-509.
$newRows = @(
    @{ Row = 506; A = "A 42531-2023"; B = 45181; C = 45184; G = 4.9 },
    @{ Row = 507; A = "A 42677-2023"; B = 45181; C = 45184; G = 0.3 },
    @{ Row = 508; A = "A 42799-2023"; B = 45182; C = 45184; G = 2.9 },
    @{ Row = 509; A = "A 42800-2023"; B = 45182; C = 45184; G = 1.2 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    $ws.Cells.Item($r, 1).Value = $rowData.A

    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = "VÄSTRA GÖTALANDS LÄN"
    $ws.Cells.Item($r, 5).Value = "ULRICEHAMN"

    $ws.Cells.Item($r, 7).Value = $rowData.G

    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($r, $col).Value = 0
    }

    $ws.Cells.Item($r, 18).Value = ""
    $ws.Cells.Item($r, 18).WrapText = $true

    # Match the source sheet's convention: every row except the very last
    # one gets an explicit 15pt row height.
    if ($r -ne 509) {
        $ws.Rows.Item($r).RowHeight = 15
    }
}
